# Babylab_WesternSydney_12-15m_bi_GF_Trial.xlsx — corrected-files upload (#41)
#
# Semantic changes applied:
#   1. Remove the empty "Sheet3" tab and rename "Sheet1" to "12-15m_Bil"
#      (the workbook now contains a single, properly-named sheet).
#   2. Fix the typo'd header label in K1: the stray, out-of-order shared
#      string "fixation_incongruent" becomes "fixation_incongrent".
#   3. Correct seven bogus "latency = 0" placeholder values (rows where the
#      trial had no valid latency) to the text marker "NA", matching the
#      convention used elsewhere in the sheet (e.g. F/G11, F/G18, ...).
#   4. Re-select G97 to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. sheet bookkeeping ---------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Delete()
$ws.Name = "12-15m_Bil"

# --- 2. header typo fix --------------------------------------------------
$ws.Range("K1").Value = "fixation_incongrent"

# --- 3. latency placeholder corrections (0 -> "NA") -----------------------
$naRows = @(10, 28, 29, 49, 63, 80, 97)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 7).Value = "NA"
}

# --- 4. cursor / selection -------------------------------------------------
$ws.Range("G97").Select()
